# Update cryptos list with latest scraped values (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = "@"
$r.Value = '64.273.30'
$r.Style = "Normal"

$r = $ws.Range('E2')
$r.NumberFormat = "@"
$r.Value = '  +0.66%  '
$r.Style = "Normal"

$r = $ws.Range('D3')
$r.NumberFormat = "@"
$r.Value = '3.489.83'
$r.Style = "Normal"

$r = $ws.Range('E3')
$r.NumberFormat = "@"
$r.Value = '  -0.16%  '
$r.Style = "Normal"

$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '586.63'
$r.Style = "Normal"

$r = $ws.Range('E5')
$r.NumberFormat = "@"
$r.Value = '  +0.32%  '
$r.Style = "Normal"

$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '134.05'
$r.Style = "Normal"

$r = $ws.Range('E6')
$r.NumberFormat = "@"
$r.Value = '  +1.65%  '
$r.Style = "Normal"

$r = $ws.Range('D7')
$r.NumberFormat = "@"
$r.Value = '3.490.28'
$r.Style = "Normal"

$r = $ws.Range('E7')
$r.NumberFormat = "@"
$r.Value = '  -0.10%  '
$r.Style = "Normal"

$r = $ws.Range('E8')
$r.NumberFormat = "@"
$r.Value = '  -0.01%  '
$r.Style = "Normal"

$r = $ws.Range('E9')
$r.NumberFormat = "@"
$r.Value = '  -0.62%  '
$r.Style = "Normal"

$r = $ws.Range('E10')
$r.NumberFormat = "@"
$r.Value = '  -0.08%  '
$r.Style = "Normal"

$r = $ws.Range('D11')
$r.NumberFormat = "@"
$r.Value = '7.18'
$r.Style = "Normal"

$r = $ws.Range('E11')
$r.NumberFormat = "@"
$r.Value = '  +1.36%  '
$r.Style = "Normal"

$r = $ws.Range('E12')
$r.NumberFormat = "@"
$r.Value = '  -2.06%  '
$r.Style = "Normal"

$r = $ws.Range('D13')
$r.NumberFormat = "@"
$r.Value = '4.087.51'
$r.Style = "Normal"

$r = $ws.Range('E13')
$r.NumberFormat = "@"
$r.Value = '  -0.31%  '
$r.Style = "Normal"

$r = $ws.Range('D14')
$r.NumberFormat = "@"
$r.Value = '0.119'
$r.Style = "Normal"

$r = $ws.Range('E14')
$r.NumberFormat = "@"
$r.Value = '  +1.87%  '
$r.Style = "Normal"

$r = $ws.Range('E15')
$r.NumberFormat = "@"
$r.Value = '  +0.20%  '
$r.Style = "Normal"

$r = $ws.Range('D16')
$r.NumberFormat = "@"
$r.Value = '3.491.00'
$r.Style = "Normal"

$r = $ws.Range('E16')
$r.NumberFormat = "@"
$r.Value = '  -0.09%  '
$r.Style = "Normal"

$r = $ws.Range('D17')
$r.NumberFormat = "@"
$r.Value = '64.317.64'
$r.Style = "Normal"

$r = $ws.Range('E17')
$r.NumberFormat = "@"
$r.Value = '  +0.59%  '
$r.Style = "Normal"

$r = $ws.Range('D18')
$r.NumberFormat = "@"
$r.Value = '25.11'
$r.Style = "Normal"

$r = $ws.Range('E18')
$r.NumberFormat = "@"
$r.Value = '  -9.35%  '
$r.Style = "Normal"

$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '9.98'
$r.Style = "Normal"

$r = $ws.Range('E19')
$r.NumberFormat = "@"
$r.Value = '  +0.06%  '
$r.Style = "Normal"

$r = $ws.Range('E20')
$r.NumberFormat = "@"
$r.Value = '  +1.65%  '
$r.Style = "Normal"

$r = $ws.Range('D21')
$r.NumberFormat = "@"
$r.Value = '13.60'
$r.Style = "Normal"

$r = $ws.Range('E21')
$r.NumberFormat = "@"
$r.Value = '  -5.58%  '
$r.Style = "Normal"

$r = $ws.Range('D22')
$r.NumberFormat = "@"
$r.Value = '387.53'
$r.Style = "Normal"

$r = $ws.Range('E22')
$r.NumberFormat = "@"
$r.Value = '  -0.72%  '
$r.Style = "Normal"

$r = $ws.Range('D23')
$r.NumberFormat = "@"
$r.Value = '3.631.40'
$r.Style = "Normal"

$r = $ws.Range('E23')
$r.NumberFormat = "@"
$r.Value = '  -0.28%  '
$r.Style = "Normal"

$r = $ws.Range('E24')
$r.NumberFormat = "@"
$r.Value = '  -2.27%  '
$r.Style = "Normal"

$r = $ws.Range('D25')
$r.NumberFormat = "@"
$r.Value = '74.59'
$r.Style = "Normal"

$r = $ws.Range('E25')
$r.NumberFormat = "@"
$r.Value = '  +2.19%  '
$r.Style = "Normal"

$r = $ws.Range('E26')
$r.NumberFormat = "@"
$r.Value = '  -0.06%  '
$r.Style = "Normal"

$r = $ws.Range('E27')
$r.NumberFormat = "@"
$r.Value = '  -0.32%  '
$r.Style = "Normal"

$r = $ws.Range('D28')
$r.NumberFormat = "@"
$r.Value = '0.0000112'
$r.Style = "Normal"

$r = $ws.Range('E28')
$r.NumberFormat = "@"
$r.Value = '  +0.04%  '
$r.Style = "Normal"

$r = $ws.Range('B29')
$r.NumberFormat = "@"
$r.Value = 'Fetch.AI'
$r.Style = "Normal"

$r = $ws.Range('C29')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$r.Style = "Normal"

$r = $ws.Range('D29')
$r.NumberFormat = "@"
$r.Value = '1.54'
$r.Style = "Normal"

$r = $ws.Range('E29')
$r.NumberFormat = "@"
$r.Value = '  -2.57%  '
$r.Style = "Normal"

$r = $ws.Range('B30')
$r.NumberFormat = "@"
$r.Value = 'Binance-PegBSC-USD'
$r.Style = "Normal"

$r = $ws.Range('C30')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$r.Style = "Normal"

$r = $ws.Range('D30')
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"

$r = $ws.Range('E30')
$r.NumberFormat = "@"
$r.Value = '  +0.30%  '
$r.Style = "Normal"

$r = $ws.Range('D31')
$r.NumberFormat = "@"
$r.Value = '7.37'
$r.Style = "Normal"

$r = $ws.Range('E31')
$r.NumberFormat = "@"
$r.Value = '  -1.54%  '
$r.Style = "Normal"

$r = $ws.Range('D32')
$r.NumberFormat = "@"
$r.Value = '8.26'
$r.Style = "Normal"

$r = $ws.Range('E32')
$r.NumberFormat = "@"
$r.Value = '  +0.16%  '
$r.Style = "Normal"

$r = $ws.Range('E33')
$r.NumberFormat = "@"
$r.Value = '  -1.08%  '
$r.Style = "Normal"

$r = $ws.Range('D34')
$r.NumberFormat = "@"
$r.Value = '3.511.06'
$r.Style = "Normal"

$r = $ws.Range('E34')
$r.NumberFormat = "@"
$r.Value = '  +0.30%  '
$r.Style = "Normal"

$r = $ws.Range('E35')
$r.NumberFormat = "@"
$r.Value = '  -0.01%  '
$r.Style = "Normal"

$r = $ws.Range('E36')
$r.NumberFormat = "@"
$r.Value = '  +2.64%  '
$r.Style = "Normal"

$r = $ws.Range('D37')
$r.NumberFormat = "@"
$r.Value = '23.43'
$r.Style = "Normal"

$r = $ws.Range('E37')
$r.NumberFormat = "@"
$r.Value = '  -1.79%  '
$r.Style = "Normal"

$r = $ws.Range('E38')
$r.NumberFormat = "@"
$r.Value = '  -1.18%  '
$r.Style = "Normal"

$r = $ws.Range('D39')
$r.NumberFormat = "@"
$r.Value = '6.82'
$r.Style = "Normal"

$r = $ws.Range('E39')
$r.NumberFormat = "@"
$r.Value = '  -2.40%  '
$r.Style = "Normal"

$r = $ws.Range('E40')
$r.NumberFormat = "@"
$r.Value = '  -2.08%  '
$r.Style = "Normal"

$r = $ws.Range('D41')
$r.NumberFormat = "@"
$r.Value = '161.53'
$r.Style = "Normal"

$r = $ws.Range('E41')
$r.NumberFormat = "@"
$r.Value = '  -3.81%  '
$r.Style = "Normal"

$r = $ws.Range('D42')
$r.NumberFormat = "@"
$r.Value = '0.0780'
$r.Style = "Normal"

$r = $ws.Range('E42')
$r.NumberFormat = "@"
$r.Value = '  -3.63%  '
$r.Style = "Normal"

$r = $ws.Range('D43')
$r.NumberFormat = "@"
$r.Value = '0.804'
$r.Style = "Normal"

$r = $ws.Range('E43')
$r.NumberFormat = "@"
$r.Value = '  -0.89%  '
$r.Style = "Normal"

$r = $ws.Range('E44')
$r.NumberFormat = "@"
$r.Value = '  -0.02%  '
$r.Style = "Normal"

$r = $ws.Range('D45')
$r.NumberFormat = "@"
$r.Value = '25.38'
$r.Style = "Normal"

$r = $ws.Range('E45')
$r.NumberFormat = "@"
$r.Value = '  -6.32%  '
$r.Style = "Normal"

$r = $ws.Range('D46')
$r.NumberFormat = "@"
$r.Value = '41.80'
$r.Style = "Normal"

$r = $ws.Range('E46')
$r.NumberFormat = "@"
$r.Value = '  -0.20%  '
$r.Style = "Normal"

$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '4.39'
$r.Style = "Normal"

$r = $ws.Range('E47')
$r.NumberFormat = "@"
$r.Value = '  +0.36%  '
$r.Style = "Normal"

$r = $ws.Range('E48')
$r.NumberFormat = "@"
$r.Value = '  -1.34%  '
$r.Style = "Normal"

$r = $ws.Range('E49')
$r.NumberFormat = "@"
$r.Value = '  +1.32%  '
$r.Style = "Normal"

$r = $ws.Range('D50')
$r.NumberFormat = "@"
$r.Value = '2.470.93'
$r.Style = "Normal"

$r = $ws.Range('E50')
$r.NumberFormat = "@"
$r.Value = '  +1.42%  '
$r.Style = "Normal"

$r = $ws.Range('D51')
$r.NumberFormat = "@"
$r.Value = '6.72'
$r.Style = "Normal"

$r = $ws.Range('E51')
$r.NumberFormat = "@"
$r.Value = '  -2.46%  '
$r.Style = "Normal"

